# ---------------------------------------------------------------------------
# 108th Precinct weekly CompStat report refresh ("New crime data collected").
#
# 1) Header text: new Police Commissioner name, report volume/number bump,
#    and the reporting week date range roll-forward.
# 2) Weekly crime-complaint grid (rows 15-28): this week's / last week's
#    counts, percent changes, and the rolling 28-day / YTD / 2-year figures
#    all got refreshed with the new week of data.
# 3) A few cells in the TOTAL/Transit/Housing block (rows 22, 31, 33) flip
#    between a numeric value and the sheet's placeholder text ("0" / "***.*")
#    used when a rate is undefined (e.g. previous count was zero).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title text -----------------------------------------------

# Police Commissioner
$ws.Range("M6").Value = "Jessica S. Tisch"

# "Volume 31   Number  47" -> "... 48"
$ws.Range("A8").Value = "Volume 31   Number  48"

# "Report Covering the Week  11/18/2024  Through  11/24/2024" -> next week
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"

# --- Weekly crime-complaint grid ----------------------------------------

# Row 15
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 100
$ws.Range("N15").Value = 75

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 218
$ws.Range("J16").Value = 240
$ws.Range("K16").Value = -9.166666666666
$ws.Range("L16").Value = 21.111111111111
$ws.Range("M16").Value = 23.863636363636
$ws.Range("N16").Value = -76.684491978609

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 242
$ws.Range("J17").Value = 221
$ws.Range("K17").Value = 9.502262443438
$ws.Range("L17").Value = 18.048780487804
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -17.123287671232

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -90
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -48
$ws.Range("I18").Value = 220
$ws.Range("J18").Value = 211
$ws.Range("K18").Value = 4.265402843601
$ws.Range("L18").Value = 38.364779874213
$ws.Range("M18").Value = -5.982905982905
$ws.Range("N18").Value = -84.647592463363

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -7.843137254901
$ws.Range("I19").Value = 663
$ws.Range("J19").Value = 685
$ws.Range("K19").Value = -3.211678832116
$ws.Range("L19").Value = 4.245283018867
$ws.Range("M19").Value = 54.545454545454
$ws.Range("N19").Value = -22.274325908558

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -10
$ws.Range("I20").Value = 237
$ws.Range("J20").Value = 285
$ws.Range("K20").Value = -16.842105263157
$ws.Range("L20").Value = 11.267605633802
$ws.Range("M20").Value = 24.736842105263
$ws.Range("N20").Value = -87.745604963805

# Row 21 (borough total)
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -38.235294117647
$ws.Range("F21").Value = 112
$ws.Range("H21").Value = -14.503816793893
$ws.Range("I21").Value = 1602
$ws.Range("J21").Value = 1666
$ws.Range("K21").Value = -3.841536614645
$ws.Range("L21").Value = 13.375796178343
$ws.Range("M21").Value = 37.275064267352
$ws.Range("N21").Value = -70.696908725077

# Row 22 (Transit) - C22 also needs special handling below (text "0" -> number)
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -28.571428571428
$ws.Range("I22").Value = 56
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = -30
$ws.Range("L22").Value = -25.333333333333
$ws.Range("M22").Value = 21.739130434782

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 59
$ws.Range("E24").Value = -32.203389830508
$ws.Range("F24").Value = 152
$ws.Range("G24").Value = 190
$ws.Range("H24").Value = -20
$ws.Range("I24").Value = 1895
$ws.Range("J24").Value = 1922
$ws.Range("K24").Value = -1.404786680541
$ws.Range("L24").Value = 34.111818825194
$ws.Range("M24").Value = 119.582850521437

# Row 25 (Retail Theft)
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = -21.875
$ws.Range("F25").Value = 100
$ws.Range("G25").Value = 111
$ws.Range("H25").Value = -9.909909909909
$ws.Range("I25").Value = 1280
$ws.Range("J25").Value = 1187
$ws.Range("K25").Value = 7.834877843302
$ws.Range("L25").Value = 121.453287197232

# Row 26 (Misd. Assault)
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 40
$ws.Range("H26").Value = -4.761904761904
$ws.Range("I26").Value = 519
$ws.Range("J26").Value = 479
$ws.Range("K26").Value = 8.350730688935
$ws.Range("L26").Value = 0.386847195357
$ws.Range("M26").Value = 12.095032397408

# Row 27 (UCR Rape*)
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 150

# Row 28 (Other Sex Crimes)
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = -42.857142857142
$ws.Range("J28").Value = 78
$ws.Range("K28").Value = -23.076923076923
$ws.Range("L28").Value = -30.232558139534

# --- Number <-> placeholder-text swaps -----------------------------------
# These cells change whether the underlying value is a real number or the
# sheet's text placeholder ("0" / "***.*") used when a rate can't be
# computed. Plain `.Value = "0"` would get auto-detected back into a number
# by the host, so instead copy an existing placeholder-text cell (same
# shared text + identical General-format style already used elsewhere in
# these rows) onto the target, which carries the text type across faithfully.

# C22: was text "0" -> now the real number 2, taking on the neighbouring
# numeric-column style (same style D22 uses).
$ws.Range("D22").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 2

# G31: was numeric 1 -> now text "0"
$ws.Range("C15").Copy($ws.Range("G31"))

# H31: was numeric -100 -> now text "***.*"
$ws.Range("E15").Copy($ws.Range("H31"))

# C33, D33: were numeric 1 -> now text "0"
$ws.Range("C15").Copy($ws.Range("C33"))
$ws.Range("D15").Copy($ws.Range("D33"))

# E33: was numeric 0 -> now text "***.*"
$ws.Range("E15").Copy($ws.Range("E33"))
